$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Facility utilisation (column F) values to match the validated
# amount for matching facility id currency and utilisation.
$ws.Range("F2").Value = 761579.37
$ws.Range("F3").Value = 761579.37
$ws.Range("F4").Value = 761579.37
$ws.Range("F5").Value = 761579.37

# Update the active selection to F2 (matches the cell last edited).
$ws.Range("F2").Select()
